# Update the "Metadata" worksheet:
#  - Experimental (row 7) value cell B7 goes from blank to the text "true"
#  - Date (row 8) value cell B8 is updated to a new timestamp
#
# Note: simply doing $ws.Range("B7").Value = "true" makes Excel
# auto-coerce the literal word "true" into a native Boolean cell,
# which is not what we want - the target keeps a plain text cell.
# Instead we build the text "true" with a formula in a scratch cell
# (so it is never parsed as a boolean literal), copy it, and paste
# only the *value* into B7. This keeps B7's original cell style/type
# (text) intact. The scratch column is then removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("Z1").Formula = '="true"'
$ws.Range("Z1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("Z1").EntireColumn.Delete()

$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
